$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.386.02"
$ws.Range("E2").Value = "  -1.81%  "

# Row 3
$ws.Range("D3").Value = "3.688.55"
$ws.Range("E3").Value = "  -2.96%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.07"
$ws.Range("E5").Value = "  -2.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.84"
$ws.Range("E6").Value = "  -3.93%  "

# Row 7
$ws.Range("D7").Value = "3.687.60"
$ws.Range("E7").Value = "  -3.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -3.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -7.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("E11").Value = "  -2.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -2.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -3.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.67"
$ws.Range("E14").Value = "  -6.03%  "

# Row 15
$ws.Range("D15").Value = "4.312.35"
$ws.Range("E15").Value = "  -2.90%  "

# Row 16
$ws.Range("D16").Value = "3.688.77"
$ws.Range("E16").Value = "  -2.52%  "

# Row 17
$ws.Range("D17").Value = "69.463.31"
$ws.Range("E17").Value = "  -1.76%  "

# Row 18
$ws.Range("E18").Value = "  -1.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.33"
$ws.Range("E19").Value = "  -5.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  -5.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "481.24"
$ws.Range("E21").Value = "  -2.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.89"
$ws.Range("E22").Value = "  -7.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.668"
$ws.Range("E23").Value = "  -7.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.31"
$ws.Range("E24").Value = "  -4.41%  "

# Row 25
$ws.Range("D25").Value = "3.834.15"
$ws.Range("E25").Value = "  -2.95%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000131"
$ws.Range("E26").Value = "  -7.94%  "

# Row 27
$ws.Range("E27").Value = "  +0.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.52"
$ws.Range("E28").Value = "  -4.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  -5.88%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.85"
$ws.Range("E30").Value = "  -8.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.74"
$ws.Range("E31").Value = "  -9.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -6.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.85"
$ws.Range("E33").Value = "  -6.07%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.168"
$ws.Range("E34").Value = "  -4.40%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.12"
$ws.Range("E35").Value = "  -6.35%  "

# Row 36
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("D37").Value = "3.654.63"
$ws.Range("E37").Value = "  -3.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.56"
$ws.Range("E38").Value = "  -5.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.30"
$ws.Range("E39").Value = "  +6.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0939"
$ws.Range("E40").Value = "  -6.67%  "

# Row 41
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("E42").Value = "  -5.23%  "

# Row 43
$ws.Range("E43").Value = "  -0.06%  "

# Row 44
$ws.Range("E44").Value = "  -6.78%  "

# Row 45
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.69"
$ws.Range("E45").Value = "  -4.03%  "

# Row 46
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -9.09%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.50"
$ws.Range("E47").Value = "  -1.23%  "

# Row 48
$ws.Range("E48").Value = "  -8.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.84"
$ws.Range("E50").Value = "  +6.09%  "

# Row 51
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "396.79"
$ws.Range("E51").Value = "  -6.38%  "
